$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.312.25"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "2.993.79"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'507.03"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").Value = "'137.95"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "'7.14"
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").Value = "'0.366"
$ws.Range("E11").Value = "  +2.36%  "
$ws.Range("D12").Value = "3.505.21"
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").Value = "'0.0000163"
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").Value = "56.236.05"
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("D17").Value = "2.992.21"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "'5.99"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "'12.92"
$ws.Range("E19").Value = "  +2.39%  "
$ws.Range("D20").Value = "'8.07"
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("D21").Value = "'331.28"
$ws.Range("E21").Value = "  +3.21%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").Value = "'0.495"
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").Value = "'64.83"
$ws.Range("E24").Value = "  +2.76%  "
$ws.Range("D25").Value = "3.116.42"
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "'0.0" + [char]0x2083 + "0943"
$ws.Range("E27").Value = "  +5.56%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.164"
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("D29").Value = "'6.35"
$ws.Range("E29").Value = "  -4.25%  "
$ws.Range("D30").Value = "'6.91"
$ws.Range("E30").Value = "  -3.22%  "
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("D32").Value = "'20.29"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("D33").Value = "'1.16"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").Value = "'153.42"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").Value = "'4.48"
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("D36").Value = "'5.83"
$ws.Range("E36").Value = "  +0.85%  "
$ws.Range("D37").Value = "'26.31"
$ws.Range("E37").Value = "  +7.62%  "
$ws.Range("D38").Value = "'1.25"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").Value = "3.026.34"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").Value = "'36.82"
$ws.Range("E41").Value = "  -2.69%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("D44").Value = "'0.651"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").Value = "2.176.06"
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").Value = "'1.35"
$ws.Range("E46").Value = "  -2.71%  "
$ws.Range("D47").Value = "'5.87"
$ws.Range("E47").Value = "  -1.37%  "
$ws.Range("D48").Value = "'0.926"
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("D50").Value = "'19.42"
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").Value = "'0.0853"
$ws.Range("E51").Value = "  -1.94%  "
